# Refresh the cryptos price/volume snapshot (GitHub Actions style update).
# Notes:
#  - Column D ("Price") values that look like plain decimal numbers
#    (e.g. "0.998", "536.96") are written with a leading apostrophe so
#    Excel stores them as literal text instead of auto-converting them
#    to numbers (matches how the sheet already stores every price as a
#    string, including thousand-dotted ones like "59.995.71").
#  - Rows 45/46 (Stellar/Hedera) swap rank order in this update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "59.995.71"
$ws.Cells.Item(2, 5).Value = "  +1.34%  "

$ws.Cells.Item(3, 4).Value = "2.656.49"
$ws.Cells.Item(3, 5).Value = "  +2.33%  "

$ws.Cells.Item(4, 5).Value = "  -0.01%  "

$ws.Cells.Item(5, 4).Value = "'536.96"
$ws.Cells.Item(5, 5).Value = "  +1.36%  "

$ws.Cells.Item(6, 4).Value = "'145.95"
$ws.Cells.Item(6, 5).Value = "  +4.00%  "

$ws.Cells.Item(7, 4).Value = "'0.998"
$ws.Cells.Item(7, 5).Value = "  +0.04%  "

$ws.Cells.Item(8, 4).Value = "'0.574"
$ws.Cells.Item(8, 5).Value = "  +1.31%  "

$ws.Cells.Item(9, 4).Value = "2.675.20"
$ws.Cells.Item(9, 5).Value = "  +2.42%  "

$ws.Cells.Item(10, 4).Value = "'6.66"
$ws.Cells.Item(10, 5).Value = "  +3.10%  "

$ws.Cells.Item(11, 5).Value = "  +2.16%  "

$ws.Cells.Item(12, 5).Value = "  +1.57%  "

$ws.Cells.Item(13, 5).Value = "  -1.24%  "

$ws.Cells.Item(14, 4).Value = "3.131.02"
$ws.Cells.Item(14, 5).Value = "  +2.43%  "

$ws.Cells.Item(15, 4).Value = "59.921.55"
$ws.Cells.Item(15, 5).Value = "  +1.34%  "

$ws.Cells.Item(16, 4).Value = "'21.34"
$ws.Cells.Item(16, 5).Value = "  +3.94%  "

$ws.Cells.Item(17, 4).Value = "2.671.81"
$ws.Cells.Item(17, 5).Value = "  +2.75%  "

$ws.Cells.Item(18, 5).Value = "  +1.27%  "

$ws.Cells.Item(19, 4).Value = "'345.28"
$ws.Cells.Item(19, 5).Value = "  -0.76%  "

$ws.Cells.Item(20, 5).Value = "  +2.10%  "

$ws.Cells.Item(21, 4).Value = "'10.30"
$ws.Cells.Item(21, 5).Value = "  +1.64%  "

$ws.Cells.Item(22, 4).Value = "'6.39"
$ws.Cells.Item(22, 5).Value = "  -0.86%  "

$ws.Cells.Item(23, 4).Value = "'0.998"
$ws.Cells.Item(23, 5).Value = "  +0.09%  "

$ws.Cells.Item(24, 4).Value = "'67.63"
$ws.Cells.Item(24, 5).Value = "  +0.42%  "

$ws.Cells.Item(25, 5).Value = "  +2.70%  "

$ws.Cells.Item(26, 5).Value = "  -0.52%  "

$ws.Cells.Item(27, 4).Value = "'0.999"
$ws.Cells.Item(27, 5).Value = "  +0.10%  "

$ws.Cells.Item(28, 5).Value = "  +2.29%  "

$ws.Cells.Item(29, 4).Value = "0.0₃0756"
$ws.Cells.Item(29, 5).Value = "  +1.98%  "

$ws.Cells.Item(30, 4).Value = "'0.998"
$ws.Cells.Item(30, 5).Value = "  -0.02%  "

$ws.Cells.Item(31, 5).Value = "  +3.02%  "

$ws.Cells.Item(32, 4).Value = "'5.92"
$ws.Cells.Item(32, 5).Value = "  +0.62%  "

$ws.Cells.Item(33, 4).Value = "'19.14"
$ws.Cells.Item(33, 5).Value = "  +1.79%  "

$ws.Cells.Item(34, 4).Value = "'150.59"
$ws.Cells.Item(34, 5).Value = "  +0.96%  "

$ws.Cells.Item(35, 4).Value = "'4.05"
$ws.Cells.Item(35, 5).Value = "  +1.53%  "

$ws.Cells.Item(36, 5).Value = "  +3.03%  "

$ws.Cells.Item(37, 5).Value = "  +0.27%  "

$ws.Cells.Item(38, 4).Value = "'0.844"
$ws.Cells.Item(38, 5).Value = "  +1.44%  "

$ws.Cells.Item(39, 4).Value = "'0.828"
$ws.Cells.Item(39, 5).Value = "  +0.61%  "

$ws.Cells.Item(40, 4).Value = "'292.69"
$ws.Cells.Item(40, 5).Value = "  +8.41%  "

$ws.Cells.Item(41, 5).Value = "  +2.20%  "

$ws.Cells.Item(42, 4).Value = "'0.998"
$ws.Cells.Item(42, 5).Value = "  +0.05%  "

$ws.Cells.Item(43, 5).Value = "  +1.51%  "

$ws.Cells.Item(44, 4).Value = "'10.74"
$ws.Cells.Item(44, 5).Value = "  -0.34%  "

$ws.Cells.Item(45, 2).Value = "Hedera"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(45, 4).Value = "'0.0539"
$ws.Cells.Item(45, 5).Value = "  +3.68%  "

$ws.Cells.Item(46, 2).Value = "Stellar"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(46, 4).Value = "'0.0957"
$ws.Cells.Item(46, 5).Value = "  -0.25%  "

$ws.Cells.Item(47, 5).Value = "  +2.39%  "

$ws.Cells.Item(48, 4).Value = "1.975.41"
$ws.Cells.Item(48, 5).Value = "  +0.63%  "

$ws.Cells.Item(49, 4).Value = "'4.61"
$ws.Cells.Item(49, 5).Value = "  -0.49%  "

$ws.Cells.Item(50, 4).Value = "'18.51"
$ws.Cells.Item(50, 5).Value = "  +1.09%  "

$ws.Cells.Item(51, 4).Value = "'109.92"
$ws.Cells.Item(51, 5).Value = "  -0.66%  "
